# Enable SELL trading: fix lots of bugs/configuration so SELL trades work.
#
# The original "Sheet1" config table had 4 data rows (rows 2-5):
#   row2: BTCUSDT / USDT   / BUY  / 1m / LMT / 100 / 100
#   row3: BTC     / Position / SELL / 1m / LMT / "Currently owned position size" / 100 (col H)
#   row4: ETHUSDT / USDT   / BUY  / 1m / LMT / 100 / 100
#   row5: ETH     / Position / SELL / 1m / LMT / "Currently owned position size" / 100 (col H)
#
# The redundant SELL helper rows (3 and 5) are removed, and the remaining
# BUY rows (now rows 2 and 3) are converted in-place into the SELL config
# that the bot actually needs (currency -> "Position", direction -> "SELL",
# the "Highest close periods" column cleared, and the position size moved
# into the "Lowest close Periods" column as 100).
#
# "Demo trading" is also flipped from FALSE to TRUE so SELL trades can be
# exercised safely first.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Remove the two now-redundant rows (delete bottom-up so row numbers
#     of the row we still need to delete don't shift) -------------------
$ws.Rows.Item(5).Delete()
$ws.Rows.Item(3).Delete()

# --- Convert the two remaining symbol rows into SELL / Position rows ---
# Row 2 (BTCUSDT)
$ws.Range("B2").Value = "Position"
$ws.Range("C2").Value = "SELL"
$ws.Range("G2").ClearContents()
$ws.Range("H2").Value = 100

# Row 3 (ETHUSDT, used to be row 4 before the deletes above)
$ws.Range("B3").Value = "Position"
$ws.Range("C3").Value = "SELL"
$ws.Range("G3").ClearContents()
$ws.Range("H3").Value = 100

# --- Flip "Demo trading" from FALSE to TRUE (now on row 6) -------------
$ws.Range("B6").Value = $true

# --- Re-create the list data validations that lived in the x14 extended
#     validation block, at their new (post row-delete) locations. The
#     standard B-column ">= 0" validation on "Open positions limit"
#     (now B5) is carried forward automatically by the row delete, but
#     the extended list validations need to be re-pointed by hand.
$ws.Range("D2:D3").Validation.Delete()
$ws.Range("D2:D3").Validation.Add(3, 1, 1, "=Validity!`$D`$2:`$D`$16")
$ws.Range("D2:D3").Validation.IgnoreBlank = $true
$ws.Range("D2:D3").Validation.InCellDropdown = $true

$ws.Range("B8").Validation.Delete()
$ws.Range("B8").Validation.Add(3, 1, 1, "=Validity!`$D`$2:`$D`$16")
$ws.Range("B8").Validation.IgnoreBlank = $true
$ws.Range("B8").Validation.InCellDropdown = $true

$ws.Range("E2:E3").Validation.Delete()
$ws.Range("E2:E3").Validation.Add(3, 1, 1, "=Validity!`$E`$2:`$E`$3")
$ws.Range("E2:E3").Validation.IgnoreBlank = $true
$ws.Range("E2:E3").Validation.InCellDropdown = $true
